# Migrating the "Guardians" link from a free-text "Ids" column (holding
# comma-separated kid ids, e.g. "1, 5, 6") to a proper numeric primary key
# (ID), and adding the matching "Guardians" foreign-key column on the Kids
# sheet (comma-separated guardian ids, e.g. "1, 2").

$wb = $excel.ActiveWorkbook

$kids = $wb.Worksheets.Item("Kids")
$guardians = $wb.Worksheets.Item("Guardians")

# --- Guardians sheet: replace the "Ids" text column with a numeric "ID" ---
$guardians.Range("A1").Value = "ID"
$guardians.Range("A2").Value = 1
$guardians.Range("A3").Value = 2
$guardians.Range("A4").Value = 3
$guardians.Range("A5").Value = 4

# --- Kids sheet: add the new "Guardians" column (E) ---
$kids.Range("E1").Value = "Guardians"
$kids.Range("E2").Value = "1, 2"
$kids.Range("E3").Value = "1, 2"
$kids.Range("E4").Value = "1, 2"
$kids.Range("E5").Value = 3
$kids.Range("E6").Value = 4

# Give the new header cell (E1) the same boxed / bold header look as the
# other header cells (copy the format from D1, the previous right-most
# header cell).
$kids.Range("D1").Copy()
$kids.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Selection / active sheet bookkeeping ---
# Final state: Kids selection sits on E7 (just below the new column) and is
# no longer the active tab; Guardians selection sits on A6 and becomes the
# active tab.
$kids.Range("E7").Select() | Out-Null
$guardians.Range("A6").Select() | Out-Null
$guardians.Activate()
